# Generate Report for Archive
# Updates the localization status from "Ready for handoff" to "In Translation"
# for the rows that had completed handoff (file-level status columns on the
# Overview sheet, and the per-file Status column on the zh-cn / de-de sheets),
# then re-fits the affected columns to their new (shorter) content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E5").Value = "In Translation"
$ws1.Range("F5").Value = "In Translation"
$ws1.Range("E6").Value = "In Translation"
$ws1.Range("F6").Value = "In Translation"
$ws1.Range("E7").Value = "In Translation"
$ws1.Range("F7").Value = "In Translation"

# --- zh-cn sheet --------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("C6").Value = "In Translation"
$ws2.Range("C7").Value = "In Translation"

# --- de-de sheet --------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("C6").Value = "In Translation"
$ws3.Range("C7").Value = "In Translation"

# --- Re-fit the Status columns now that the text is shorter -------------
$ws1.Columns.Item(5).AutoFit() | Out-Null
$ws1.Columns.Item(6).AutoFit() | Out-Null
$ws2.Columns.Item(3).AutoFit() | Out-Null
$ws3.Columns.Item(3).AutoFit() | Out-Null
